$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells (shorten column titles, fix typos)
# Order matters for how new shared strings get appended (matches author's edit order)
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Formula"
$ws.Range("A1").Value = "ID"
$ws.Range("D1").Value = "Mass"

# Match the author's final selection
$ws.Range("H12").Select()
